$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 621 (Excel shifts rows 621:644 down to 622:645,
# and copies formatting from the row above into the new row, matching column D's date style).
$ws.Rows.Item(621).Insert()

# Populate the newly inserted row 621 with the new weekly record.
$ws.Range("A621").Value = 3
$ws.Range("B621").Value = "Femacal de La Calera"
$ws.Range("C621").Value = "Coquimbo"
$ws.Range("D621").Value = 45075
$ws.Range("E621").Value = 5
$ws.Range("F621").Value = 100112017
$ws.Range("G621").Value = "Apio"
$ws.Range("H621").Value = "Americana (o)"
$ws.Range("I621").Value = "Primera"
$ws.Range("J621").Value = 335
$ws.Range("K621").Value = 9000
$ws.Range("L621").Value = 9500
$ws.Range("M621").Value = 9246
$ws.Range("N621").Value = "$/docena de matas"
$ws.Range("O621").Value = "Provincia de Limarí"
$ws.Range("P621").Value = 1541
$ws.Range("Q621").Value = 6
$ws.Range("R621").Value = "Hortaliza"
